$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 3 (col C) text + reflow rows 2-5 existing content with updated alignment ---

# Row 1
$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("A1").VerticalAlignment = -4160
$ws.Range("B1").Value = "Test cases"
$ws.Range("B1").HorizontalAlignment = -4131
$ws.Range("B1").VerticalAlignment = -4160
$ws.Range("C1").Value = "Test steps"
$ws.Range("C1").HorizontalAlignment = -4131
$ws.Range("C1").VerticalAlignment = -4160
$ws.Range("D1").Value = "Test Data"
$ws.Range("D1").HorizontalAlignment = -4131
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("E1").Value = "Expected Result"
$ws.Range("E1").HorizontalAlignment = -4131
$ws.Range("E1").VerticalAlignment = -4160

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("A2").HorizontalAlignment = -4131
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("B2").Value = "Verify item in Mobile List page can be sorted by 'Name'"
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B2").VerticalAlignment = -4160
$ws.Range("C2").Value = "1. Goto http://live.techoanda.org/ `n2. Verify Title of the page`n3. Click on 'MOBILE' menu`n4. Verify Title of the page `n5. In the list of all mobile , select 'SORT BY' dropdown as 'name'`n6. Verify all productssorted by name"
$ws.Range("C2").HorizontalAlignment = -4131
$ws.Range("C2").VerticalAlignment = -4160
$ws.Range("C2").WrapText = $true
$ws.Range("D2").HorizontalAlignment = -4131
$ws.Range("D2").VerticalAlignment = -4160
$ws.Range("E2").Value = "1. Text 'THIS IS DEMO SITE FOR   ' shown in home page`n2. Title 'MOBILE' is shown on mobile list page.`n3.  All 3 products sorted by name"
$ws.Range("E2").HorizontalAlignment = -4131
$ws.Range("E2").VerticalAlignment = -4160
$ws.Range("E2").WrapText = $true
$ws.Range("A2:E2").EntireRow.RowHeight = 90

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("B3").Value = "Verify that cost of product in list page and details page are equal"
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").VerticalAlignment = -4160
$ws.Range("C3").Value = "1.Goto http://live.techpanda.org/`n2. Click on 'MOBILE' menu`n3. In the list of all mobile , read the cost of Sony Xperia mobile. Note this value `n4.Click on Sony Xperia mobile `n5. Read the Sony Xperia mobile from detail page.`n6. Compare value in Step 3 & Step 5"
$ws.Range("C3").HorizontalAlignment = -4131
$ws.Range("C3").VerticalAlignment = -4160
$ws.Range("C3").WrapText = $true
$ws.Range("D3").HorizontalAlignment = -4131
$ws.Range("D3").VerticalAlignment = -4160
$ws.Range("E3").Value = "1) Product Value in list and details page should be equal (`$IOO)"
$ws.Range("E3").HorizontalAlignment = -4131
$ws.Range("E3").VerticalAlignment = -4160
$ws.Range("A3:E3").EntireRow.RowHeight = 90

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4160
$ws.Range("B4").Value = "Verify that you cannot add more product in cart than the product available in store"
$ws.Range("B4").HorizontalAlignment = -4131
$ws.Range("B4").VerticalAlignment = -4160
$ws.Range("C4").Value = "1.Goto http://live.techpanda.org/`n2. Click on 'MOBILE' menu`n3. In the list of all mobile , click on 'ADD TO CART' for Sony Xperia mobile `n4. Change 'QTY' value to 1000 and click 'UPDATE' button. `n5. Verify the error message `n6. Then click on 'EMPTY CART' link in the footer of list of all mobiles. `n7. Verify cart is empty"
$ws.Range("C4").HorizontalAlignment = -4131
$ws.Range("C4").VerticalAlignment = -4160
$ws.Range("C4").WrapText = $true
$ws.Range("D4").HorizontalAlignment = -4131
$ws.Range("D4").VerticalAlignment = -4160
$ws.Range("E4").Value = "1. On clicking update button an error is shown 'The requested quantity for Sony Xperia`" is not available.'`n2. On clicking empty cart button - a message 'SHOPPING CART IS EMPTY' is shown"
$ws.Range("E4").HorizontalAlignment = -4131
$ws.Range("E4").VerticalAlignment = -4160
$ws.Range("E4").WrapText = $true
$ws.Range("A4:E4").EntireRow.RowHeight = 105

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4160
$ws.Range("B5").Value = "Verify that you are able to compare two product"
$ws.Range("B5").HorizontalAlignment = -4131
$ws.Range("B5").VerticalAlignment = -4160
$ws.Range("C5").Value = "1. Goto http://live.techpanda.org/ `n2.Click on 'MOBILE' menu `n3.In mobile products list , click on 'Add To Compare' for 2 mobiles`n 4 Click on 'COMPARE' button.`n5. Verify the pop-up window and check that the products are reflected in it `n6 Close the Popup Windows"
$ws.Range("C5").HorizontalAlignment = -4131
$ws.Range("C5").VerticalAlignment = -4160
$ws.Range("C5").WrapText = $true
$ws.Range("D5").Value = "Phone 1 Sony Xperia`nPhone 2 IPHONE"
$ws.Range("D5").HorizontalAlignment = -4131
$ws.Range("D5").VerticalAlignment = -4160
$ws.Range("D5").WrapText = $true
$ws.Range("E5").Value = "1.A Popup window opens with headind as 'COMPARE PRODUCTS' and the selected products are present in it. `n2. Popup window is closed "
$ws.Range("E5").HorizontalAlignment = -4131
$ws.Range("E5").VerticalAlignment = -4160
$ws.Range("E5").WrapText = $true
$ws.Range("A5:E5").EntireRow.RowHeight = 90

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("A6").HorizontalAlignment = -4131
$ws.Range("A6").VerticalAlignment = -4160
$ws.Range("B6").Value = "Verify you can create account in E-commerce site and can share wishlist to other people using email"
$ws.Range("B6").HorizontalAlignment = -4131
$ws.Range("B6").VerticalAlignment = -4160
$ws.Range("C6").Value = "1.Goto http://live.techpanda.org/ `n2.Click on my account link`n 3.Click Create Account link and fill New User information except Email ID`n 4.Click Register`n 5. Verify Registration is done `n6.Go to TV menu `n7.Add product in your wish list `n8.Click SHARE WISHLIST `n9.1n next page enter Email and a message and click SHARE WISHLIST `n10. Check wishlist is shared"
$ws.Range("C6").HorizontalAlignment = -4131
$ws.Range("C6").VerticalAlignment = -4160
$ws.Range("C6").WrapText = $true
$ws.Range("D6").Value = "product = LG LCD"
$ws.Range("D6").HorizontalAlignment = -4131
$ws.Range("D6").VerticalAlignment = -4160
$ws.Range("D6").WrapText = $true
$ws.Range("E6").Value = "1. Account Registration done `n2. Wishlist Shared Successfully"
$ws.Range("E6").HorizontalAlignment = -4131
$ws.Range("E6").VerticalAlignment = -4160
$ws.Range("E6").WrapText = $true
$ws.Range("A6:E6").EntireRow.RowHeight = 150

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("A7").HorizontalAlignment = -4131
$ws.Range("A7").VerticalAlignment = -4160
$ws.Range("B7").Value = "Verify user is able to purchase Continue product using registered email id(USE CHROME BROWSER)"
$ws.Range("B7").HorizontalAlignment = -4131
$ws.Range("B7").VerticalAlignment = -4160
$ws.Range("C7").Value = "1 Go to http://live.techpanda.org/ `n2.Click on my account link `n3.Login in application using previously created credential `n4.Click on MY WISHLIST link `n5.1n next page, Click ADD TO CART link`n 6.Click PROCEED TO CHECKOUT`n 7. Enter Shipping Information`n 8. Click Estimate `n9.Verify Shipping cost generated`n10. Select Shipping Cost, Update Total `n11. Verify shipping cost is added to total `n12Click `"Proceed to Checkout`n13. Enter Billing Information `n14. In Shipping Method, Click Continue `n15.1n Payment Information select 'Check/Money Order' radio button.Click `n16.C1ick 'PLACE ORDER' button`n17. Verify Oder is generated. Note the order number"
$ws.Range("C7").HorizontalAlignment = -4131
$ws.Range("C7").VerticalAlignment = -4160
$ws.Range("C7").WrapText = $true
$ws.Range("D7").Value = "User credentials created in test case 05`n Shipping Information Country United States `nState = New York`n Zip = 542896`n Address ABC `nCity = New York `nTelephone = 12345678"
$ws.Range("D7").HorizontalAlignment = -4131
$ws.Range("D7").VerticalAlignment = -4160
$ws.Range("D7").WrapText = $true
$ws.Range("E7").Value = "1. Flat Rate Shipping of `$5 is generated `n2. Shipping cost is added to total product cost`n3. Oder is placed. Order number is generated"
$ws.Range("E7").Characters(1, 91).Font.Bold = $false
$ws.Range("E7").Characters(92, 2).Font.Bold = $true
$ws.Range("E7").Characters(94, 41).Font.Bold = $false
$ws.Range("E7").HorizontalAlignment = -4131
$ws.Range("E7").VerticalAlignment = -4160
$ws.Range("E7").WrapText = $true
$ws.Range("A7:E7").EntireRow.RowHeight = 285

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("A8").HorizontalAlignment = -4131
$ws.Range("A8").VerticalAlignment = -4160
$ws.Range("B8").Value = "Verify that you will be able to save previously placed order as a pdf file"
$ws.Range("B8").HorizontalAlignment = -4131
$ws.Range("B8").VerticalAlignment = -4160
$ws.Range("C8").Value = "1.Go to http://live.techpanda.org/ `n2.Click on My Account link `n3.Login in application using previously created credential `n4.Click on 'My Orders' `n5.Click on 'View Order `n6.Verify the previously created order is displayed in 'RECENT ORDERS' table and status is Pending `n7.Click on 'Print Order link `n8. Verify Order is saved as PDF"
$ws.Range("C8").HorizontalAlignment = -4131
$ws.Range("C8").VerticalAlignment = -4160
$ws.Range("C8").WrapText = $true
$ws.Range("D8").Value = "Use FireFox"
$ws.Range("D8").HorizontalAlignment = -4131
$ws.Range("D8").VerticalAlignment = -4160
$ws.Range("E8").Value = "Previously created order is displayed in 'RECENT ORDERS' table and status is Pending. 2`n.Order is saved as PDF."
$ws.Range("E8").HorizontalAlignment = -4131
$ws.Range("E8").VerticalAlignment = -4160
$ws.Range("E8").WrapText = $true
$ws.Range("A8:E8").EntireRow.RowHeight = 120

# --- View state: selection + zoom (best-effort; window chrome itself is host-managed) ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 115
$ws.Range("C4").Select()

